# Edit script: update weather-prediction literature review sheet
#  - normalize stray font colors/sizes back to the workbook's default theme font
#  - add newly-researched rows of data (Dataset used / Platform Used / Methods used)
#  - update two Focus/Aim descriptions with fuller text
#  - resize a couple of rows/columns to fit the new content

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Reset font on the whole populated data block (A1:K6) back to the
#    workbook default (Calibri 11, automatic/theme colour). This clears the
#    assorted FF333333 / FF000000 / Arial-14 overrides that had accumulated
#    from copy-pasting content from the web.
# ---------------------------------------------------------------------------
$dataRange = $ws.Range("A1:K6")
$dataRange.Font.Name = "Calibri"
$dataRange.Font.Size = 11
$dataRange.Font.ColorIndex = -4105   # xlColorIndexAutomatic

# ---------------------------------------------------------------------------
# 2. Update existing Focus/Aim text (column F) for rows 2 and 4.
# ---------------------------------------------------------------------------
$ws.Range("F2").Value = "The focus is to survey the state-of-the-art studies of deep learning-based weather forecasting, in the aspects of the design of neural network (NN) architectures, spatial and temporal scales, as well as the datasets and benchmarks. Then we analyze the advantages and disadvantages of DLWP by comparing it with the conventional NWP, and summarize the potential future research topics of DLWP."

$ws.Range("F4").Value = "Predicting the weather by analyzing features like temperature,apparent temperature,humidity,wind speed, wind bearing, visibility, cloud cover with Random Forest, Decision Tree, MLP classifier, Linear regression, and Gaussian naive Bayes are examples of machine learning methods."

# ---------------------------------------------------------------------------
# 3. Fill in the newly-added "Dataset used" / "Platform Used" / "Methods used"
#    cells for each reference row.
# ---------------------------------------------------------------------------
$era5 = "ERA5 (European Centre for Medium-Range Weather Forecasts (ECMWF) Reanalysis 5th Generation)"

$ws.Range("G2").Value = $era5
$ws.Range("H2").Value = "The artificial neural network(ANN) "

$ws.Range("G3").Value = $era5
$h3 = $ws.Range("H3")
$h3.Value = "Huawei Atlas AI infrastructure"
$h3.Characters(16, 15).Font.Name = "Calibri"
$h3.Characters(16, 15).Font.Size = 11
$h3.Characters(16, 15).Font.ColorIndex = -4105
$ws.Range("I3").Value = "Pangu-Weather model"

$ws.Range("G4").Value = "Kaggle.com"
$ws.Range("I4").Value = "Random Forest, Decision Tree, MLP classifier, Linear regression, and Gaussian naive Bayes "

$ws.Range("G5").Value = $era5
$ws.Range("H5").Value = "The Hadoop Distributed File System"

$ws.Range("I6").Value = "ARIMA model"

# ---------------------------------------------------------------------------
# 4. Alignment touch-ups: several cells gained WrapText / centred horizontal
#    alignment once content was normalised.
# ---------------------------------------------------------------------------
$ws.Range("A2:A6").WrapText = $true
$ws.Range("B2:B6").WrapText = $true
$ws.Range("D2").WrapText = $true
$ws.Range("D4").WrapText = $true
$ws.Range("E2:E6").WrapText = $true

$centerWrapCells = @("G2","H2","G3","H3","I3","G4","I4","G5","H5","I6")
foreach ($addr in $centerWrapCells) {
    $c = $ws.Range($addr)
    $c.HorizontalAlignment = -4108   # xlCenter
    $c.VerticalAlignment = -4108     # xlCenter
    $c.WrapText = $true
}

$ws.Range("H4").VerticalAlignment = -4108
$ws.Range("H4").WrapText = $true

# ---------------------------------------------------------------------------
# 5. Row heights / column width adjustments to fit the new text.
# ---------------------------------------------------------------------------
$ws.Rows.Item(3).RowHeight = 126
$ws.Rows.Item(4).RowHeight = 120
$ws.Rows.Item(5).RowHeight = 172.9
$ws.Rows.Item(6).RowHeight = 75

$ws.Columns.Item(9).ColumnWidth = 17.29

# ---------------------------------------------------------------------------
# 6. Selection / view bookkeeping to mirror the saved workbook state.
# ---------------------------------------------------------------------------
$ws.Range("I6").Select()
